$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Status (C3) becomes "Success"
$ws.Range("C3").Value = "Success"

# Row 4: Pass (B4) value changes and becomes left-aligned
$ws.Range("B4").Value = 12233
$ws.Range("B4").HorizontalAlignment = -4131  # xlLeft
